$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10, columns A-T
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myl9"
$ws.Range("C2").Value = "Cd69"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 10.953471
$ws.Range("H2").Value = 32.860413
$ws.Range("I2").Value = 0.1262416469639394
$ws.Range("J2").Value = 0.1262416469639394
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 13.92410533333333
$ws.Range("N2").Value = 41.772316
$ws.Range("O2").Value = 0.9926481335887254
$ws.Range("P2").Value = 0.9926481335887255
$ws.Range("Q2").Value = 152.517283969612
$ws.Range("R2").Value = 1372.655555726508
$ws.Range("S2").Value = 0.1253135352399213
$ws.Range("T2").Value = 0.1253135352399213

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myl9"
$ws.Range("C3").Value = "Cd69"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 10.953471
$ws.Range("H3").Value = 32.860413
$ws.Range("I3").Value = 0.1262416469639394
$ws.Range("J3").Value = 0.1262416469639394
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05023433333333333
$ws.Range("N3").Value = 0.150703
$ws.Range("O3").Value = 0.003581200804767964
$ws.Range("P3").Value = 0.003581200804767964
$ws.Range("Q3").Value = 0.550240313371
$ws.Range("R3").Value = 4.952162820339001
$ws.Range("S3").Value = 0.0004520966877024931
$ws.Range("T3").Value = 0.0004520966877024931

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myl9"
$ws.Range("C4").Value = "Cd69"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 10.953471
$ws.Range("H4").Value = 32.860413
$ws.Range("I4").Value = 0.1262416469639394
$ws.Range("J4").Value = 0.1262416469639394
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.052892
$ws.Range("N4").Value = 0.158676
$ws.Range("O4").Value = 0.003770665606506583
$ws.Range("P4").Value = 0.003770665606506583
$ws.Range("Q4").Value = 0.579350988132
$ws.Range("R4").Value = 5.214158893188
$ws.Range("S4").Value = 0.0004760150363156725
$ws.Range("T4").Value = 0.0004760150363156725

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myl9"
$ws.Range("C5").Value = "Cd69"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.453103
$ws.Range("H5").Value = 13.359309
$ws.Range("I5").Value = 0.05132318849614514
$ws.Range("J5").Value = 0.05132318849614514
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 13.92410533333333
$ws.Range("N5").Value = 41.772316
$ws.Range("O5").Value = 0.9926481335887254
$ws.Range("P5").Value = 0.9926481335887255
$ws.Range("Q5").Value = 62.00547523218265
$ws.Range("R5").Value = 558.049277089644
$ws.Range("S5").Value = 0.05094586727052082
$ws.Range("T5").Value = 0.05094586727052083

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myl9"
$ws.Range("C6").Value = "Cd69"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.453103
$ws.Range("H6").Value = 13.359309
$ws.Range("I6").Value = 0.05132318849614514
$ws.Range("J6").Value = 0.05132318849614514
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05023433333333333
$ws.Range("N6").Value = 0.150703
$ws.Range("O6").Value = 0.003581200804767964
$ws.Range("P6").Value = 0.003581200804767964
$ws.Range("Q6").Value = 0.2236986604696666
$ws.Range("R6").Value = 2.013287944227
$ws.Range("S6").Value = 0.0001837986439456529
$ws.Range("T6").Value = 0.0001837986439456529

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myl9"
$ws.Range("C7").Value = "Cd69"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.453103
$ws.Range("H7").Value = 13.359309
$ws.Range("I7").Value = 0.05132318849614514
$ws.Range("J7").Value = 0.05132318849614514
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.052892
$ws.Range("N7").Value = 0.158676
$ws.Range("O7").Value = 0.003770665606506583
$ws.Range("P7").Value = 0.003770665606506583
$ws.Range("Q7").Value = 0.235533523876
$ws.Range("R7").Value = 2.119801714884
$ws.Range("S7").Value = 0.0001935225816786688
$ws.Range("T7").Value = 0.0001935225816786688

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Myl9"
$ws.Range("C8").Value = "Cd69"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 71.35933300000001
$ws.Range("H8").Value = 214.077999
$ws.Range("I8").Value = 0.8224351645399155
$ws.Range("J8").Value = 0.8224351645399154
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 13.92410533333333
$ws.Range("N8").Value = 41.772316
$ws.Range("O8").Value = 0.9926481335887254
$ws.Range("P8").Value = 0.9926481335887255
$ws.Range("Q8").Value = 993.6148692084093
$ws.Range("R8").Value = 8942.533822875683
$ws.Range("S8").Value = 0.8163887310782834
$ws.Range("T8").Value = 0.8163887310782834

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Myl9"
$ws.Range("C9").Value = "Cd69"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 71.35933300000001
$ws.Range("H9").Value = 214.077999
$ws.Range("I9").Value = 0.8224351645399155
$ws.Range("J9").Value = 0.8224351645399154
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05023433333333333
$ws.Range("N9").Value = 0.150703
$ws.Range("O9").Value = 0.003581200804767964
$ws.Range("P9").Value = 0.003581200804767964
$ws.Range("Q9").Value = 3.584688520366333
$ws.Range("R9").Value = 32.262196683297
$ws.Range("S9").Value = 0.002945305473119819
$ws.Range("T9").Value = 0.002945305473119818

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Myl9"
$ws.Range("C10").Value = "Cd69"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 71.35933300000001
$ws.Range("H10").Value = 214.077999
$ws.Range("I10").Value = 0.8224351645399155
$ws.Range("J10").Value = 0.8224351645399154
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.052892
$ws.Range("N10").Value = 0.158676
$ws.Range("O10").Value = 0.003770665606506583
$ws.Range("P10").Value = 0.003770665606506583
$ws.Range("Q10").Value = 3.774337841036
$ws.Range("R10").Value = 33.969040569324
$ws.Range("S10").Value = 0.003101127988512242
$ws.Range("T10").Value = 0.003101127988512241
